$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain text even when the value looks like a plain number
    # (e.g. "308.29"), then restore the default "Normal" style so we do not leave
    # a stray number-format behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.845.34"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "1.815.67"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("E5").Value = "  -0.34%  "
Set-TextValue $ws.Range("D6") "308.29"
$ws.Range("E6").Value = "  -1.90%  "
Set-TextValue $ws.Range("D7") "0.4609"
$ws.Range("E7").Value = "  -2.73%  "
Set-TextValue $ws.Range("D8") "0.3637"
$ws.Range("E8").Value = "  -1.45%  "
Set-TextValue $ws.Range("D9") "0.07212"
$ws.Range("E9").Value = "  -3.27%  "
Set-TextValue $ws.Range("D10") "0.8571"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.878.74"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.07500"
$ws.Range("E13").Value = "  +2.07%  "
Set-TextValue $ws.Range("D14") "5.314"
$ws.Range("E14").Value = "  -2.45%  "
Set-TextValue $ws.Range("D15") "6.493"
$ws.Range("E15").Value = "  -1.31%  "
Set-TextValue $ws.Range("D16") "91.42"
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("E17").Value = "  -0.29%  "
Set-TextValue $ws.Range("D18") "0.000008566"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D20") "14.39"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "26.716.50"
$ws.Range("E21").Value = "  -3.22%  "
Set-TextValue $ws.Range("D22") "5.126"
$ws.Range("E22").Value = "  -3.53%  "
Set-TextValue $ws.Range("D23") "10.48"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("D24").Value = "1.983.85"
$ws.Range("E24").Value = "  -5.63%  "
Set-TextValue $ws.Range("D25") "151.36"
$ws.Range("E25").Value = "  -0.49%  "
Set-TextValue $ws.Range("D26") "1.845"
$ws.Range("E26").Value = "  -3.39%  "
Set-TextValue $ws.Range("D27") "18.03"
$ws.Range("E27").Value = "  -3.25%  "
Set-TextValue $ws.Range("D28") "2.071"
$ws.Range("E28").Value = "  -3.35%  "
Set-TextValue $ws.Range("D29") "5.082"
$ws.Range("E29").Value = "  -3.09%  "
Set-TextValue $ws.Range("D30") "115.17"
$ws.Range("E30").Value = "  -2.12%  "
Set-TextValue $ws.Range("D31") "0.08845"
$ws.Range("E31").Value = "  -1.68%  "
Set-TextValue $ws.Range("D32") "2.949"
Set-TextValue $ws.Range("D33") "4.404"
$ws.Range("E33").Value = "  -3.19%  "
Set-TextValue $ws.Range("D34") "1.127"
$ws.Range("E34").Value = "  -4.34%  "
Set-TextValue $ws.Range("D35") "0.7111"
$ws.Range("E35").Value = "  -5.98%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D36") "1.072"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D37") "0.05231"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D38") "2.404"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.01912"
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.919"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "7.116"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D42") "0.5126"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D43") "0.1618"
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D44") "8.135"
$ws.Range("E44").Value = "  -4.43%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.4776"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D46") "1.008"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D47") "10.07"
$ws.Range("E47").Value = "  -4.87%  "
Set-TextValue $ws.Range("D48") "102.91"
$ws.Range("E48").Value = "  -1.90%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.06251"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "1.612"
$ws.Range("E50").Value = "  -3.80%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "63.83"
$ws.Range("E51").Value = "  -2.99%  "
